$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row: "_old" columns become "_FV2310", "_new" columns become "_FV2404".
$headersFV2310 = @("Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310")
for ($i = 0; $i -lt $headersFV2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2310[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

$headersFV2404 = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2404[$i]
}

# Turn the data range into an Excel Table ("ListObject") with the new header names.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
